$wb = $excel.ActiveWorkbook

# Header values (row 1) shared by every results sheet in this workbook.
$headers = @(5, 50, 500, 5000, 50000, 500000)

# Row labels (column A) shared by every results sheet.
$labels = @("byte", "int", "string", "date")

# New sheet: same_elements
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sameElements = $wb.Worksheets.Add($null, $lastSheet)
$sameElements.Name = "same_elements"

for ($c = 0; $c -lt 6; $c++) {
    $sameElements.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$sameElementsData = @(
    @(0, 0, 0.001001, 0.029029, 0.273287, 3.377133),
    @(0, 0, 0.001501, 0.019021, 0.241287, 3.05628),
    @(0, 0, 0.001001, 0.01952, 0.246221, 3.44815),
    @(0, 0.000501, 0.001502, 0.017017, 0.265277, 3.747413)
)

for ($i = 0; $i -lt 4; $i++) {
    $row = $i + 2
    $sameElements.Cells.Item($row, 1).Value = $labels[$i]
    for ($c = 0; $c -lt 6; $c++) {
        $sameElements.Cells.Item($row, $c + 2).Value = $sameElementsData[$i][$c]
    }
}

# New sheet: partly_same
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$partlySame = $wb.Worksheets.Add($null, $lastSheet)
$partlySame.Name = "partly_same"

for ($c = 0; $c -lt 6; $c++) {
    $partlySame.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$partlySameData = @(
    @(0, 0, 0.001001, 0.024525, 0.299313, 3.352512),
    @(0, 0.000501, 0.001502, 0.027029, 0.289802, 4.007186),
    @(0, 0, 0.006507, 0.032533, 0.403922, 5.936201),
    @(0, 0.0005, 0.002002, 0.029031, 0.359913, 4.692449)
)

for ($i = 0; $i -lt 4; $i++) {
    $row = $i + 2
    $partlySame.Cells.Item($row, 1).Value = $labels[$i]
    for ($c = 0; $c -lt 6; $c++) {
        $partlySame.Cells.Item($row, $c + 2).Value = $partlySameData[$i][$c]
    }
}

# "same_elements" is the active tab after this edit.
$sameElements.Select()
